# seedDB gender to lower case
#
# The roster's "phone" column (F) was never populated for the last 13
# people added to the sheet (rows 42-54). Fill it in with the next block
# of sequential seed phone numbers, continuing on from the last one
# already present in F41 (612-616-3287).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$phones = @(
    "612-616-3288",
    "612-616-3289",
    "612-616-3290",
    "612-616-3291",
    "612-616-3292",
    "612-616-3293",
    "612-616-3294",
    "612-616-3295",
    "612-616-3296",
    "612-616-3297",
    "612-616-3298",
    "612-616-3299",
    "612-616-3300"
)

$firstRow = 42
for ($i = 0; $i -lt $phones.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 6).Value = $phones[$i]
}

# Mirror the author's view state: scrolled down so row 38 is at the top,
# with the newly-filled phone cells selected.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 38
$win.ScrollColumn = 1
$ws.Range("F41:F54").Select()
